$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "68.794.38"
$ws.Range("E2").Value = "  +2.02%  "
$ws.Range("D3").Value = "3.762.24"
$ws.Range("E3").Value = "  -0.67%  "
$ws.Range("E4").Value = "  -0.22%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "603.32"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.80%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "169.52"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.26%  "
$ws.Range("D7").Value = "3.759.20"
$ws.Range("E7").Value = "  -0.62%  "
$ws.Range("E8").Value = "  -0.05%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.536"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +1.84%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.164"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +4.36%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.34"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +2.43%  "
$ws.Range("E12").Value = "  +0.42%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "38.44"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.67%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.0000246"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.58%  "
$ws.Range("D15").Value = "4.389.34"
$ws.Range("E15").Value = "  -0.82%  "
$ws.Range("D16").Value = "3.759.45"
$ws.Range("E16").Value = "  -1.05%  "
$ws.Range("D17").Value = "68.803.16"
$ws.Range("E17").Value = "  +1.87%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "7.33"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +2.47%  "
$ws.Range("E19").Value = "  +0.18%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "17.23"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.27%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "10.79"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +15.65%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "497.28"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.25%  "
$ws.Range("E23").Value = "  -0.56%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "85.81"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.70%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.0000146"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.17%  "
$ws.Range("E26").Value = "  -1.11%  "
$ws.Range("E27").Value = "  +1.00%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.27"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.94%  "
$ws.Range("E29").Value = "  -0.09%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.54"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +5.19%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.00"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +1.04%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "7.99"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +2.83%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "32.21"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.45%  "
$ws.Range("D34").Value = "3.908.65"
$ws.Range("E34").Value = "  -0.76%  "
$ws.Range("D35").Value = "3.694.75"
$ws.Range("E35").Value = "  -0.89%  "
$ws.Range("E36").Value = "  +0.15%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.999"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.27%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.01"
$ws.Range("D38").Style = "Normal"
$ws.Range("E39").Value = "  +0.78%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.132"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.87%  "
$ws.Range("E41").Value = "  -0.06%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "444.56"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -4.20%  "
$ws.Range("E43").Value = "  -1.02%  "
$ws.Range("E44").Value = "  +0.05%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.87"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +1.20%  "
$ws.Range("E46").Value = "  +1.89%  "
$ws.Range("E47").Value = "  +0.01%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "40.73"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.06%  "
$ws.Range("D49").Value = "2.838.86"
$ws.Range("E49").Value = "  +0.19%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "142.09"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +1.36%  "
$ws.Range("E51").Value = "  +2.45%  "
